$d = $word.ActiveDocument

# --- Change 1 & 2: consolidate split runs in the first table's 1st cell ---
# "a" + " cell " -> "a cell "
$d.Content.Find.Execute("a" + [char]0x0B + " cell ", $true, $false, $false, $false, $false, $true, 1, $false, "a cell ", 2) | Out-Null
$d.Content.Find.Execute("a cell ", $true, $false, $false, $false, $false, $true, 1, $false, "a cell ", 2) | Out-Null

# "w" + "ith three paragraphs" + " with auto margin" + "s" + "." -> "with three paragraphs containing auto margins."
$d.Content.Find.Execute("with three paragraphs with auto margins.", $true, $false, $false, $false, $false, $true, 1, $false, "with three paragraphs containing auto margins.", 2) | Out-Null

# --- Change 3: add a new row to the first table ---
$t = $d.Tables.Item(1)
$t.Rows.Add() | Out-Null
$cell = $t.Cell($t.Rows.Count, 1)
$r = $cell.Range

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:spacing w:beforeAutospacing="0"/><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t>This is</w:t></w:r></w:p>
<w:p><w:pPr><w:spacing w:beforeAutospacing="0"/><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t>a cell </w:t></w:r></w:p>
<w:p><w:pPr><w:spacing w:beforeAutospacing="0"/><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="en-US"/></w:rPr><w:t>with three paragraphs with disabled auto margins.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xml)

Write-Host "edit complete"
